$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 28-63: mark as "yes" and fill in the feature/java file references
for ($r = 28; $r -le 63; $r++) {
    $ws.Range("B$r").Value = "yes"
}

for ($r = 28; $r -le 35; $r++) {
    $ws.Range("E$r").Value = "HTB_straight_flush_beats_AIP.feature"
    $ws.Range("F$r").Value = "HTBStraightFlushBeatsAIP.java"
}
for ($r = 36; $r -le 42; $r++) {
    $ws.Range("E$r").Value = "HTB_four_of_a_kind_beats_AIP.feature"
    $ws.Range("F$r").Value = "HTBFourOfAKindBeatsAIP.java"
}
for ($r = 43; $r -le 48; $r++) {
    $ws.Range("E$r").Value = "HTB_full_house_beats_AIP.feature"
    $ws.Range("F$r").Value = "HTBFullHouseBeatsAIP.java"
}
for ($r = 49; $r -le 53; $r++) {
    $ws.Range("E$r").Value = "HTB_flush_beats_AIP.feature"
    $ws.Range("F$r").Value = "HTBFlushBeatsAIP.java"
}
for ($r = 54; $r -le 57; $r++) {
    $ws.Range("E$r").Value = "HTB_straight_beats_AIP.feature"
    $ws.Range("F$r").Value = "HTBStraightBeatsAIP.java"
}
for ($r = 58; $r -le 60; $r++) {
    $ws.Range("E$r").Value = "HTB_three_of_a_kind_beats_AIP.feature"
    $ws.Range("F$r").Value = "HTBThreeOfAKindBeatsAIP.java"
}
for ($r = 61; $r -le 62; $r++) {
    $ws.Range("E$r").Value = "HTB_two_pair_beats_AIP.feature"
    $ws.Range("F$r").Value = "HTBTwoPairBeatsAIP.java"
}
$ws.Range("E63").Value = "HTB_pair_beats_AIP.feature"
$ws.Range("F63").Value = "HTBPairBeatsAIP.java"

# Update the active selection to match the saved view state
$ws.Range("C63").Select()
